$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark: remove it from its old spot (right
#    before "Start with 17w12" in the Product Version bullet) - it will be
#    re-created at the site of the new edit further down (matches Word's
#    own behaviour of moving _GoBack to the location of the last edit).
# ---------------------------------------------------------------------------
$rGoBackOld = $d.Range(0, 0)
$rGoBackOld.Find.Execute("Start with 17w12", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rGoBackOld.Collapse(1)
$rGoBackOld.MoveStart(1, -2)   # include the manual line break + the "S" of "Start"
$savedText = $rGoBackOld.Text
$rGoBackOld.Delete()
$rGoBackOld.InsertAfter($savedText)

# ---------------------------------------------------------------------------
# 2. Instance Types bullet: "Proposal: use i3" -> "Use i3"
# ---------------------------------------------------------------------------
$r1 = $d.Range(0, 0)
$r1.Find.Execute("Proposal: use i3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Text = "Use i3"

# ---------------------------------------------------------------------------
# 3. Instance Types bullet: add " by default" right after
#    "... instance type family" (before the trailing ". ").
# ---------------------------------------------------------------------------
$r2 = $d.Range(0, 0)
$r2.Find.Execute("instance type family", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)
$r2.InsertAfter(" by default")

# ---------------------------------------------------------------------------
# 4. Merge the "Instance Types" bullet with the "Compliance" bullet,
#    inserting a new line "Allow R3 and R4 families as options (less and
#    slower storage)?" plus a blank line, and re-create the "_GoBack"
#    bookmark (zero length) at that new edit location.
# ---------------------------------------------------------------------------
$r3 = $d.Range(0, 0)
$r3.Find.Execute(" restart.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Collapse(0)
$r3.MoveEnd(1, 1)     # select the paragraph mark ending the Instance Types bullet
$r3.Delete()           # merge the two paragraphs into one

$r3.InsertAfter([char]11 + "Allow R3 and R4 families as options (less and slower storage)?" + [char]11)

# Put the (recreated) "_GoBack" bookmark between the new sentence and the
# trailing line break, i.e. right before the very last inserted character.
$bmPos = $d.Range($r3.End - 1, $r3.End - 1)
$d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null

# ---------------------------------------------------------------------------
# 5. Footer page number field cached text: "1" -> "2"
# ---------------------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2) | Out-Null
